$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update China totals (row 4) ---
$ws.Range("B4").Value = 80967
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = 71150
$ws.Range("E4").Value = 6569
$ws.Range("F4").Value = 2136
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 3248

# --- Update Estados Unidos totals (row 9) ---
$ws.Range("B9").Value = 13920
$ws.Range("C9").Value = 4661
$ws.Range("D9").Value = 121
$ws.Range("E9").Value = 13585
$ws.Range("F9").Value = 64
$ws.Range("G9").Value = 64
$ws.Range("H9").Value = 214

# --- Nueva Zelanda moves up the ranking (was row 103) and gets updated
#     figures; Tunez..Liechtenstein each shift down one row (96-103) ---
$ws.Range("A96").Value = "Nueva Zelanda"
$ws.Range("B96").Value = 39
$ws.Range("C96").Value = 19
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 39
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0

$ws.Range("A97").Value = "Tunez"
$ws.Range("B97").Value = 39
$ws.Range("C97").Value = 10
$ws.Range("D97").Value = 1
$ws.Range("E97").Value = 37
$ws.Range("F97").Value = 2
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 1

$ws.Range("A98").Value = "Camboya"
$ws.Range("B98").Value = 37
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 1
$ws.Range("E98").Value = 36
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0

$ws.Range("A99").Value = "Senegal"
$ws.Range("B99").Value = 36
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 2
$ws.Range("E99").Value = 34
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0

$ws.Range("A100").Value = "Republica Dominicana"
$ws.Range("B100").Value = 34
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 32
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 2

$ws.Range("A101").Value = "Guadalupe"
$ws.Range("B101").Value = 33
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 33
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0

$ws.Range("A102").Value = "Burkina Faso"
$ws.Range("B102").Value = 33
$ws.Range("C102").Value = 6
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = 32
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 1

$ws.Range("A103").Value = "Liechtenstein"
$ws.Range("B103").Value = 28
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 0
$ws.Range("E103").Value = 28
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 0
